# fileCls imports all spreadsheets
#
# Adds a second worksheet ("Another sample") that mirrors the existing
# "Sample" sheet (same layout/formatting, imported as its own sheet), and
# records a new data point (Temperature, column C) on row 7 of the
# original "Sample" sheet.

$wb = $excel.ActiveWorkbook

# The original worksheet that gets "imported" a second time.
$ws1 = $wb.Worksheets.Item("Sample")

# Duplicate it (preserves formatting/merged cells/styles) and place the
# copy immediately after the original, then give it its own title.
$ws1.Copy($null, $ws1)
$ws2 = $wb.Worksheets.Item(2)
$ws2.Name = "Another sample"
$ws2.Range("A1").Value = "Sample spreadsheet for fpdt testing, second sheet"

# New reading recorded on the original sheet: column C ("Temperature 3")
# for day 5 (row 7). Copy the formatting from the neighbouring cell so the
# new cell matches the rest of the data column, then set its value.
$ws1.Range("A7").Copy()
$ws1.Range("C7").PasteSpecial(-4122)
$ws1.Range("C7").Value = 175.0

# Leave the original sheet selected/active, as it was before the edit.
$ws1.Activate() | Out-Null
$ws1.Range("A1").Select() | Out-Null
